# questions_points.xlsx - "new exams and bugs"
# Insert a new exam file (ttest-08.Rnw, 8 points) into the list and
# correct the point values that had drifted for the existing ttest rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the new exam file "ttest-08.Rnw" right before anova-01.Rnw (row 26)
$ws.Rows("26:26").Insert()
$ws.Range("A26").Value = "ttest-08.Rnw"
$ws.Range("B26").Value = 8

# Correct the point values for the existing ttest rows (bug fix)
$ws.Range("B20").Value = 8
$ws.Range("B21").Value = 12
$ws.Range("B22").Value = 11
$ws.Range("B23").Value = 10
$ws.Range("B24").Value = 8
$ws.Range("B25").Value = 9

# Keep the conditional formatting ("<7" highlight) covering the full points column
$ws.Range("B2:B81").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("B2:B82"))

# Restore the cursor/selection position as left by the author
$ws.Range("C21").Select()
